# Updates cryptos list data (prices and 1h volume changes) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.263.36'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '''1.606.06'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''212.62'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = '''0.486'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +0.59%  '
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D10').Value = '''18.46'
$ws.Range('E10').Value = '  +2.77%  '
$ws.Range('D11').Value = '''0.0814'
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('D12').Value = '''1.830.64'
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('D13').Value = '''1.612.44'
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('D15').Value = '''0.515'
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('D16').Value = '''26.244.23'
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('E17').Value = '  +2.76%  '
$ws.Range('D18').Value = '''0.0₃0728'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').Value = '''200.54'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('D23').Value = '''6.02'
$ws.Range('E23').Value = '  +0.60%  '
$ws.Range('D24').Value = '''1.86'
$ws.Range('E24').Value = '  +2.27%  '
$ws.Range('D25').Value = '''143.96'
$ws.Range('E25').Value = '  +1.71%  '
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('D28').Value = '''15.21'
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('E29').Value = '  +2.33%  '
$ws.Range('D30').Value = '''0.0494'
$ws.Range('E30').Value = '  +5.06%  '
$ws.Range('D31').Value = '''1.18'
$ws.Range('E31').Value = '  +0.57%  '
$ws.Range('E32').Value = '  +2.98%  '
$ws.Range('E33').Value = '  -1.13%  '
$ws.Range('E34').Value = '  +1.32%  '
$ws.Range('E35').Value = '  +1.08%  '
$ws.Range('D36').Value = '''1.163.76'
$ws.Range('E36').Value = '  +3.94%  '
$ws.Range('E37').Value = '  +3.20%  '
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('E41').Value = '  +1.38%  '
$ws.Range('D42').Value = '''5.36'
$ws.Range('E42').Value = '  +4.57%  '
$ws.Range('E43').Value = '  +0.72%  '
$ws.Range('D44').Value = '''1.741.42'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').Value = '''92.20'
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '''1.54'
$ws.Range('E46').Value = '  +1.91%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '''0.0₆0105'
$ws.Range('E47').Value = '  +13.41%  '
$ws.Range('D48').Value = '''54.09'
$ws.Range('E48').Value = '  +1.23%  '
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('E51').Value = '  -0.20%  '
